# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" sheets, reflecting a newer scrape/generation run.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1157
$wsExhibition.Range("F11").Value = 2295
$wsExhibition.Range("F13").Value = 1308
$wsExhibition.Range("F17").Value = 739
$wsExhibition.Range("F24").Value = 4459
$wsExhibition.Range("F32").Value = 648
$wsExhibition.Range("F40").Value = 89
$wsExhibition.Range("F41").Value = 124
$wsExhibition.Range("F42").Value = 110

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F5").Value = 1157
$wsAllTypes.Range("F17").Value = 2295
$wsAllTypes.Range("F19").Value = 1308
$wsAllTypes.Range("F24").Value = 739
$wsAllTypes.Range("F29").Value = 4459
$wsAllTypes.Range("F37").Value = 648
$wsAllTypes.Range("F44").Value = 89
$wsAllTypes.Range("F45").Value = 124
$wsAllTypes.Range("F46").Value = 110
